$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (above current row 5),
# pushing the existing data rows (5-14) down to rows (7-16). Excel carries
# the row-5 formatting (date style on column D) onto the newly inserted rows.
$ws.Rows("5:6").Insert()

# Populate the two newly inserted rows with the latest weekly observations.
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 44487
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112026
$ws.Range("G5").Value = "Haba"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 8000
$ws.Range("N5").Value = "`$/saco 25 kilos"
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 320
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44487
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112026
$ws.Range("G6").Value = "Haba"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 9000
$ws.Range("N6").Value = "`$/saco 25 kilos"
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 360
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
